$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price column (D) from Excel auto-coercing numeric-looking
# text (e.g. "1.00", "4.80") into real numbers: force Text format while
# writing, then clear the formatting back to General so the saved style
# matches the original (unstyled) cells but the value stays literal text.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "59.663.72"
$ws.Range("E2").Value = "  +8.17%  "
$ws.Range("D3").Value = "2.579.56"
$ws.Range("E3").Value = "  +9.54%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "504.58"
$ws.Range("E5").Value = "  +5.77%  "
$ws.Range("D6").Value = "156.33"
$ws.Range("E6").Value = "  +7.52%  "
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  +2.46%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "2.574.42"
$ws.Range("E9").Value = "  +9.64%  "
$ws.Range("E10").Value = "  +13.44%  "
$ws.Range("E11").Value = "  +5.59%  "
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +5.38%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "3.006.15"
$ws.Range("E14").Value = "  +9.05%  "
$ws.Range("D15").Value = "59.215.46"
$ws.Range("E15").Value = "  +7.64%  "
$ws.Range("D16").Value = "21.85"
$ws.Range("E16").Value = "  +8.85%  "
$ws.Range("E17").Value = "  +4.79%  "
$ws.Range("D18").Value = "2.576.69"
$ws.Range("E18").Value = "  +9.55%  "
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("D20").Value = "335.02"
$ws.Range("E20").Value = "  +5.68%  "
$ws.Range("D21").Value = "10.35"
$ws.Range("E21").Value = "  +7.48%  "
$ws.Range("E22").Value = "  +7.02%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "59.69"
$ws.Range("E24").Value = "  +5.04%  "
$ws.Range("E25").Value = "  +5.31%  "
$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  +7.29%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "2.647.33"
$ws.Range("E28").Value = "  +8.07%  "
$ws.Range("D29").Value = "0.0₃0830"
$ws.Range("E29").Value = "  +9.64%  "
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "157.29"
$ws.Range("E32").Value = "  +6.73%  "
$ws.Range("E33").Value = "  +6.51%  "
$ws.Range("E34").Value = "  +5.91%  "
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +8.25%  "
$ws.Range("E36").Value = "  +9.57%  "
$ws.Range("D37").Value = "3.91"
$ws.Range("E37").Value = "  +8.49%  "
$ws.Range("D38").Value = "0.846"
$ws.Range("E38").Value = "  +3.63%  "
$ws.Range("E39").Value = "  +9.90%  "
$ws.Range("E40").Value = "  +7.54%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "293.95"
$ws.Range("E41").Value = "  +16.83%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "35.10"
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").Value = "0.626"
$ws.Range("E44").Value = "  +8.67%  "
$ws.Range("D45").Value = "0.0567"
$ws.Range("E45").Value = "  +8.69%  "
$ws.Range("D46").Value = "0.775"
$ws.Range("E46").Value = "  +23.52%  "
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "19.16"
$ws.Range("E48").Value = "  +14.52%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0236"
$ws.Range("E49").Value = "  +6.24%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.80"
$ws.Range("E50").Value = "  +10.16%  "
$ws.Range("D51").Value = "2.002.29"
$ws.Range("E51").Value = "  +12.13%  "

$priceRange.ClearFormats()

Write-Output "done"
